# Update countries & provincias Spain
#
# The source feed was re-pulled at a later timestamp (01:52 -> 02:22), which
# shifted several countries' rank in the (descending, by total-cases) table
# and refreshed the day's counters. This reproduces that refreshed snapshot
# by writing the new timestamp, the updated counters, and the country labels
# in their new row positions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range('A1').Value = 'Datos actualizados a 5 de Abril de 2020 a las 02:22'

# Estados Unidos (row 4) - refreshed counters
$ws.Range("B4").Value = 310286
$ws.Range("C4").Value = 33125
$ws.Range("E4").Value = 287093
$ws.Range("G4").Value = 1048
$ws.Range("H4").Value = 8452

# Rows 32-36: India overtakes Malasia/Ecuador/Japon/Filipinas in the ranking
$ws.Range('A32').Value = 'India'
$ws.Range("B32").Value = 3588
$ws.Range("C32").Value = 529
$ws.Range("D32").Value = 229
$ws.Range("E32").Value = 3273
$ws.Range("F32").Value = 0
$ws.Range("G32").Value = 0
$ws.Range("H32").Value = 86

$ws.Range('A33').Value = 'Malasia'
$ws.Range("B33").Value = 3483
$ws.Range("C33").Value = 150
$ws.Range("D33").Value = 915
$ws.Range("E33").Value = 2511
$ws.Range("F33").Value = 99
$ws.Range("G33").Value = 4
$ws.Range("H33").Value = 57

$ws.Range('A34').Value = 'Ecuador'
$ws.Range("B34").Value = 3465
$ws.Range("C34").Value = 97
$ws.Range("D34").Value = 100
$ws.Range("E34").Value = 3193
$ws.Range("F34").Value = 100
$ws.Range("G34").Value = 27
$ws.Range("H34").Value = 172

$ws.Range('A35').Value = 'Japon'
$ws.Range("B35").Value = 3139
$ws.Range("C35").Value = 204
$ws.Range("D35").Value = 514
$ws.Range("E35").Value = 2548
$ws.Range("F35").Value = 64
$ws.Range("H35").Value = 77

$ws.Range('A36').Value = 'Filipinas'
$ws.Range("B36").Value = 3094
$ws.Range("C36").Value = 76
$ws.Range("D36").Value = 57
$ws.Range("E36").Value = 2893
$ws.Range("F36").Value = 1
$ws.Range("G36").Value = 8
$ws.Range("H36").Value = 144

# Republica de Chipre (row 86) - refreshed counters
$ws.Range("E86").Value = 384
$ws.Range("G86").Value = 2
$ws.Range("H86").Value = 9

# Rows 144-145: Zambia overtakes Puerto Rico
$ws.Range('A144').Value = 'Zambia'
$ws.Range("D144").Value = 2
$ws.Range("H144").Value = 1

$ws.Range('A145').Value = 'Puerto Rico'
$ws.Range("D145").Value = 1
$ws.Range("H145").Value = 2

# Rows 147-148: Bermudas overtakes Islas Caimanes
$ws.Range('A147').Value = 'Bermudas'
$ws.Range("B147").Value = 37
$ws.Range("C147").Value = 2
$ws.Range("D147").Value = 14
$ws.Range("E147").Value = 23
$ws.Range("H147").Value = 0

$ws.Range('A148').Value = 'Islas Caimanes'
$ws.Range("C148").Value = 6
$ws.Range("D148").Value = 1
$ws.Range("E148").Value = 33
$ws.Range("H148").Value = 1

# Rows 177-178: Seychelles overtakes Laos (counters tied, only labels swap)
$ws.Range('A177').Value = 'Seychelles'

$ws.Range('A178').Value = 'Laos'

# Rows 181-183: Sudan/Liberia/Angola reshuffle
$ws.Range('A181').Value = 'Sudan'
$ws.Range("C181").Value = 0
$ws.Range("D181").Value = 2
$ws.Range("G181").Value = 0
$ws.Range("H181").Value = 2

$ws.Range('A182').Value = 'Liberia'
$ws.Range("C182").Value = 3
$ws.Range("D182").Value = 3
$ws.Range("G182").Value = 1
$ws.Range("H182").Value = 1

$ws.Range('A183').Value = 'Angola'
$ws.Range("C183").Value = 2

# Rows 184 & 186: Republica del Chad overtakes San Cristobal y Nieves
$ws.Range('A184').Value = 'Republica del Chad'
$ws.Range("C184").Value = 1

$ws.Range('A186').Value = 'San Cristobal y Nieves'
$ws.Range("C186").Value = 0

# Rows 187-188: Zimbabue overtakes Nepal
$ws.Range('A187').Value = 'Zimbabue'
$ws.Range("C187").Value = 0
$ws.Range("D187").Value = 0
$ws.Range("H187").Value = 1

$ws.Range('A188').Value = 'Nepal'
$ws.Range("C188").Value = 3
$ws.Range("D188").Value = 1
$ws.Range("H188").Value = 0

# Rows 192-194: Somalia/San Vicente y las Granadinas/Cabo Verde reshuffle
$ws.Range('A192').Value = 'Somalia'
$ws.Range("C192").Value = 0
$ws.Range("D192").Value = 1
$ws.Range("H192").Value = 0

$ws.Range('A193').Value = 'San Vicente y las Granadinas'

$ws.Range('A194').Value = 'Cabo Verde'
$ws.Range("C194").Value = 1
$ws.Range("D194").Value = 0
$ws.Range("H194").Value = 1

# Rows 200-203: Belice/Sierra Leona/Malaui/Sahara Occidental reshuffle
$ws.Range('A200').Value = 'Belice'
$ws.Range("C200").Value = 0

$ws.Range('A201').Value = 'Sierra Leona'
$ws.Range("C201").Value = 2

$ws.Range('A202').Value = 'Malaui'
$ws.Range("C202").Value = 1

$ws.Range('A203').Value = 'Sahara Occidental'
$ws.Range("C203").Value = 4
